$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: new entry for "Mehr Arbeit an Kollision" / "8h" / date 2018-10-28
$ws.Range("B23").Value = "Mehr Arbeit an Kollision"
$ws.Range("D23").Value = "8h"
$ws.Range("F23").NumberFormat = "d-mmm"
$d1 = Get-Date -Year 2018 -Month 10 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Range("F23").Value = $d1

# Row 25: new entry for "Mehr Arbeit an Kollision…" / "10h" / date 2018-10-29
$ws.Range("B25").Value = "Mehr Arbeit an Kollision…"
$ws.Range("D25").Value = "10h"
$ws.Range("F25").NumberFormat = "d-mmm"
$d2 = Get-Date -Year 2018 -Month 10 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("F25").Value = $d2

# Update view: scroll so row 7 is at top, and select B27 as active cell
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("B27").Select()
